$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $c = $ws.Range($cell)
    $c.Value = "'" + $value
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '66.116.49'
$ws.Range('E2').Value = '  -2.06%  '
$ws.Range('D3').Value = '2.496.19'
$ws.Range('E3').Value = '  -5.11%  '
Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.04%  '
Set-TextValue 'D5' '579.27'
$ws.Range('E5').Value = '  -2.61%  '
Set-TextValue 'D6' '170.01'
$ws.Range('E6').Value = '  +0.75%  '
$ws.Range('E7').Value = '  +0.03%  '
Set-TextValue 'D8' '0.519'
$ws.Range('E8').Value = '  -2.89%  '
$ws.Range('D9').Value = '2.494.69'
$ws.Range('E9').Value = '  -5.14%  '
$ws.Range('E10').Value = '  -2.93%  '
$ws.Range('E11').Value = '  -0.44%  '
Set-TextValue 'D12' '0.348'
$ws.Range('E12').Value = '  -4.83%  '
Set-TextValue 'D13' '5.08'
$ws.Range('E13').Value = '  -2.74%  '
Set-TextValue 'D14' '26.34'
$ws.Range('E14').Value = '  -5.00%  '
$ws.Range('E15').Value = '  -5.46%  '
Set-TextValue 'D16' '0.0000174'
$ws.Range('E16').Value = '  -4.88%  '
$ws.Range('D17').Value = '65.716.65'
$ws.Range('E17').Value = '  -2.53%  '
$ws.Range('D18').Value = '2.496.38'
$ws.Range('E18').Value = '  -5.01%  '
Set-TextValue 'D19' '11.16'
$ws.Range('E19').Value = '  -7.38%  '
Set-TextValue 'D20' '7.63'
$ws.Range('E20').Value = '  -5.20%  '
Set-TextValue 'D21' '343.52'
$ws.Range('E21').Value = '  -3.75%  '
Set-TextValue 'D22' '4.17'
$ws.Range('E22').Value = '  -3.54%  '
Set-TextValue 'D23' '4.54'
$ws.Range('E23').Value = '  -3.05%  '
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('E25').Value = '  -0.68%  '
Set-TextValue 'D26' '69.03'
$ws.Range('E26').Value = '  -0.79%  '
Set-TextValue 'D27' '9.79'
$ws.Range('E27').Value = '  -5.43%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').Value = '2.619.89'
$ws.Range('E29').Value = '  -5.15%  '
$ws.Range('D30').Value = '0.0₃0959'
$ws.Range('E30').Value = '  -4.89%  '
Set-TextValue 'D31' '523.89'
$ws.Range('E31').Value = '  -4.22%  '
Set-TextValue 'D32' '8.04'
$ws.Range('E32').Value = '  +1.15%  '
$ws.Range('E33').Value = '  -3.07%  '
$ws.Range('E34').Value = '  -4.58%  '
Set-TextValue 'D35' '0.131'
$ws.Range('E35').Value = '  -3.69%  '
Set-TextValue 'D36' '0.999'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D37' '1.44'
$ws.Range('E37').Value = '  -4.26%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D38' '156.22'
$ws.Range('E38').Value = '  -0.12%  '
Set-TextValue 'D39' '18.45'
$ws.Range('E39').Value = '  -3.12%  '
Set-TextValue 'D40' '18.32'
$ws.Range('E40').Value = '  +0.64%  '
Set-TextValue 'D41' '0.352'
$ws.Range('E41').Value = '  -4.12%  '
$ws.Range('E42').Value = '  -3.09%  '
Set-TextValue 'D43' '5.04'
$ws.Range('E43').Value = '  -3.60%  '
$ws.Range('E44').Value = '  -0.07%  '
Set-TextValue 'D45' '2.42'
$ws.Range('E45').Value = '  -0.47%  '
Set-TextValue 'D46' '146.19'
$ws.Range('E46').Value = '  -4.40%  '
Set-TextValue 'D47' '0.553'
$ws.Range('E47').Value = '  -4.86%  '
Set-TextValue 'D48' '3.65'
$ws.Range('E48').Value = '  -3.84%  '
Set-TextValue 'D49' '1.71'
$ws.Range('E49').Value = '  +0.71%  '
$ws.Range('E50').Value = '  -9.86%  '
Set-TextValue 'D51' '0.0749'
$ws.Range('E51').Value = '  -2.83%  '
